# Apply the Jogos_da_Semana_FlashScore_2024-11-23.xlsx update.
#
# Summary of the change (per the supplied diff):
#  - The match that used to be on row 4 (Brazil - Serie A Betano, Sao Paulo
#    vs Atletico-MG) was removed from the sheet.
#  - The match that used to be on row 5 (Mexico - Liga de Expansion MX,
#    Tapatio vs Celaya) moved up to become row 4, and its
#    Odd_Under05_FT (column N) value was updated from 11 to 11.5.
#  - A handful of odds on the remaining rows were refreshed:
#      M2 (Odd_Over05_FT):  1.1  -> 1.13
#      O2 (Odd_Over15_FT):  1.58 -> 1.62
#      BD2 (Odd_CS_4-4_HT): 151  -> 126
#      Q3 (Odd_Over25_FT):  2.05 -> 2
#  - The used range / dimension shrinks from A1:BD5 to A1:BD4 as a result
#    of the row removal.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old row 4; this shifts the old row 5 up to become the new
# row 4 (matching the diff, which effectively drops the Brazil match and
# keeps the Mexico match, renumbered from 5 to 4).
$ws.Rows(4).Delete()

# Odd value refreshes on row 2.
$ws.Range("M2").Value = 1.13
$ws.Range("O2").Value = 1.62
$ws.Range("BD2").Value = 126

# Odd value refresh on row 3.
$ws.Range("Q3").Value = 2

# The Mexico match (now row 4) had its Odd_Under05_FT value updated.
$ws.Range("N4").Value = 11.5
